# Generate Report for Handoff
# Updates the "411e37d7-efb5-4e7c-bf8c-4fd2b121f893.md" row with fresh
# handoff / xliff-generation timestamps across the Overview, zh-cn and
# de-de sheets (a new handoff xliff was generated for that file).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-26 04:40:17"

# --- zh-cn sheet: column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-26 04:40:12"

# --- de-de sheet: column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-26 04:40:17"
